$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.054.32'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.122.28'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '531.47'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.27'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +4.59%  '
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.409'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.79%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '3.664.42'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.53'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000163'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '58.103.40'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.126.22'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.98'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.62'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.09'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '353.78'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.77'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.503'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0₃0883'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.26'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.13'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '21.29'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.00'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.15'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '158.60'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.06'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '25.96'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.25'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0671'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.699'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '37.59'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.54%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.399.61'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.54%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.166.44'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.977'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.03'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '19.82'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.739'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0906'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.00%  '
